# Loan RBI, Variable Instalments
# Insert a new (blank) column N on the "Repayment schedule" sheet, shifting
# the existing Late / Date / Outstanding columns one place to the right,
# and make the "Repayment schedule" sheet the active tab/selection.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

# Insert a blank column before column N (14), shifting N:P -> O:Q
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = 9.86

# Transactions sheet keeps its own selection, but is no longer the active tab
$wsTrans.Select()
$wsTrans.Range("D3").Select()

# Repayment schedule becomes the active tab/selected sheet
$wsRepay.Select()
$wsRepay.Range("R7").Select()
